$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F40").Value = 1800
$ws.Range("F41").Value = 1600
$ws.Range("F42").Value = 1680

$ws.Range("F49").Select()
